$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove all existing hyperlinks (and their formatting) so the sheet
#        can be rebuilt cleanly; we re-add the correct ones further down. ---
$ws.Cells.Hyperlinks.Delete()
$ws.Range("A1:P2").Style = "Normal"

# --- 2. Row 1: column headers -------------------------------------------
$ws.Range("A1").Value = "タイトル"
$ws.Range("B1").Value = "説明"
$ws.Range("C1").Value = "ライセンス"
$ws.Range("D1").Value = "アイテムURL"
$ws.Range("E1").Value = "dcterms:references"
$ws.Range("F1").Value = "機械可読ドキュメント"
$ws.Range("G1").Value = "dcterms:isReferencedBy"
$ws.Range("H1").Value = "年"
$ws.Range("I1").Value = "IIIFマニフェストURI"
$ws.Range("J1").Value = "viewingDirection"
$ws.Range("K1").Value = "帰属"
$ws.Range("L1").Value = "ID"
$ws.Range("M1").Value = "ソート用項目"
$ws.Range("N1").Value = "コレクション"
$ws.Range("O1").Value = "サムネイル"
$ws.Range("P1").Value = "ウェブサイトURL"

# --- 3. Row 2: values ------------------------------------------------------
$ws.Range("A2").Value = "松乃栄"
# B2 keeps its existing long description text untouched.
$ws.Range("C2").Value = "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse"
$ws.Range("D2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/document/c416f868-754f-4fed-9974-6ba911e2c0ba"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/21824"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/21824/manifest"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "東京大学総合図書館 General Library in the University of Tokyo, JAPAN"
$ws.Range("L2").Value = "c416f868-754f-4fed-9974-6ba911e2c0ba"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "松乃栄"
$ws.Range("O2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/21827/full/200,151/0/default.jpg"
$ws.Range("P2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/"

# --- 4. Re-create hyperlinks on the cells that hold URLs -------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/document/c416f868-754f-4fed-9974-6ba911e2c0ba") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/21824") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/21824/manifest") | Out-Null
$ws.Hyperlinks.Add($ws.Range("O2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/21827/full/200,151/0/default.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("P2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/matsunosakae/") | Out-Null
